$wb = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item("测试总况")
$wsDetail = $wb.Worksheets.Item("测试详情")

# ------------------------------------------------------------------
# Sheet "测试总况" (overview) updates
# ------------------------------------------------------------------
# Test run timestamp (row 6, column B) and elapsed time (row 6, column D).
$wsOverview.Range("B6").Value = "2018-03-22 17:17:09"
$wsOverview.Range("D6").Value = "193秒"

# Pass counts bumped from 2 to 3 (total cases row, pass-count row, and the
# per-device pass-count row at the bottom of the sheet).
$wsOverview.Range("D3").Value = 3
$wsOverview.Range("D4").Value = 3
$wsOverview.Range("B9").Value = 3

# ------------------------------------------------------------------
# Sheet "测试详情" (detail) updates
# ------------------------------------------------------------------
# Row 3 ("testLogin" case): precondition / steps / checkpoint text updated.
$wsDetail.Range("E3").Value = "数据不清空，已登出"
$wsDetail.Range("F3").Value = "输入用户名`n输入密码`n点击登陆`n"
$wsDetail.Range("G3").Value = "美容顾问姓名`n"

# Row 4 previously described a duplicate "testLogin2" case; it now becomes a
# brand-new "testMyExclusive" test case. The checkpoint column (G) held the
# same text as row 3's, so it must be refreshed to match row 3's new text too.
$wsDetail.Range("C4").Value = "浏览我的专属"
$wsDetail.Range("D4").Value = "testMyExclusive"
$wsDetail.Range("E4").Value = "已经登陆并停留在主页"
$wsDetail.Range("F4").Value = "打开我的订单`n取消设定促销时间`n打开我的专属`n打开全部订单`n打开成员订单状态`n返回上一级菜单`n打开全部订单`n返回上一级菜单`n打开最新订单`n返回上一级菜单`n返回主页`n"
$wsDetail.Range("G4").Value = "美容顾问姓名`n"

# Row 5 was entirely empty; copy row 4's formatting down (keeps the border /
# centered-alignment style used throughout the table) and then fill in the
# new "testOpenOrderingUI" test case values.
$wsDetail.Range("A4:J4").Copy($wsDetail.Range("A5:J5"))
$wsDetail.Range("C5").Value = "打开我的订单"
$wsDetail.Range("D5").Value = "testOpenOrderingUI"
$wsDetail.Range("F5").Value = "打开我的订单`n取消设定促销时间`n返回主页`n"
$wsDetail.Range("G5").Value = "美容顾问姓名`n"
